$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3, shifting rows 3-5 down to 4-6
$ws.Rows.Item(3).Insert()

# Populate the newly inserted row 3 with data (copy of old row 3 with updated Fecha and Volumen)
$ws.Cells.Item(3, 1).Value = 4
$ws.Cells.Item(3, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(3, 3).Value = "Los Lagos"
$ws.Cells.Item(3, 4).Value = 44469
$ws.Cells.Item(3, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(3, 5).Value = 10
$ws.Cells.Item(3, 6).Value = 100112012
$ws.Cells.Item(3, 7).Value = "Espinaca"
$ws.Cells.Item(3, 8).Value = "Sin especificar"
$ws.Cells.Item(3, 9).Value = "Primera"
$ws.Cells.Item(3, 10).Value = 20
$ws.Cells.Item(3, 11).Value = 12000
$ws.Cells.Item(3, 12).Value = 12000
$ws.Cells.Item(3, 13).Value = 12000
$ws.Cells.Item(3, 14).Value = "$/cuna 10 kilos"
$ws.Cells.Item(3, 15).Value = "Región Metropolitana"
$ws.Cells.Item(3, 16).Value = 1200
$ws.Cells.Item(3, 17).Value = 10
$ws.Cells.Item(3, 18).Value = "Hortaliza"
